# Gantt chart update from 5/1/2020
#
# Adds a new worksheet "200501" (a snapshot of the Gantt chart as of
# 5/1/2020) at the end of the workbook, alongside the existing "200413"
# and "200422" sheets, and updates the previously-active sheet's
# selection/tab state now that the new sheet is the active one.

$wb = $excel.ActiveWorkbook

function Set-Cell($sheet, $row, $col, $value, $Style) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.Value = $value
    switch ($Style) {
        1 { $cell.NumberFormat = "m/d/yy" }       # matches existing style s="1" (numFmtId 14)
        2 { $cell.NumberFormat = "mm/dd/yy;@" }    # matches existing style s="2" (numFmtId 164)
        3 { $cell.WrapText = $true }               # matches existing style s="3" (wrapText)
        default { }
    }
}

$ws2 = $wb.Worksheets.Item("200422")
$ws1 = $wb.Worksheets.Item("200413")

# Add the new sheet at the end of the workbook and name it "200501".
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "200501"

# Date values below are written as raw day-count serials (matching how
# Excel itself stores dates in the sheet XML); NumberFormat controls how
# they are displayed. The trailing numeric argument selects a style:
# 0 = none, 1 = "m/d/yy" (short date), 2 = "mm/dd/yy;@", 3 = wrap text.
Set-Cell $ws 1 1 "Task" 0
Set-Cell $ws 1 2 "Category" 0
Set-Cell $ws 1 3 "Start Date" 0
Set-Cell $ws 1 4 "End Date" 0
Set-Cell $ws 2 1 "Descriptor/Property correlation" 0
Set-Cell $ws 2 2 "Initial data" 0
Set-Cell $ws 2 3 43934 2  # 2020-04-13
Set-Cell $ws 2 4 43943 2  # 2020-04-22
Set-Cell $ws 3 1 "ML exploration (NN, RFR, KRR)" 0
Set-Cell $ws 3 2 "Initial data" 0
Set-Cell $ws 3 3 43936 2  # 2020-04-15
Set-Cell $ws 3 4 43943 2  # 2020-04-22
Set-Cell $ws 4 1 "Train models with DFT data" 0
Set-Cell $ws 4 2 "Initial data" 0
Set-Cell $ws 4 3 43936 2  # 2020-04-15
Set-Cell $ws 4 4 43943 2  # 2020-04-22
Set-Cell $ws 5 1 "Clean DFT data: outliers, normalize" 3
Set-Cell $ws 5 2 "Total data" 0
Set-Cell $ws 5 3 43943 2  # 2020-04-22
Set-Cell $ws 5 4 43950 2  # 2020-04-29
Set-Cell $ws 6 1 "Descriptor importance: remove unecessary" 3
Set-Cell $ws 6 2 "Total data" 0
Set-Cell $ws 6 3 43945 2  # 2020-04-24
Set-Cell $ws 6 4 43950 2  # 2020-04-29
Set-Cell $ws 7 1 "Train models with new DFT data (NN, RFR, LASSO, GPR)" 3
Set-Cell $ws 7 2 "Total data" 0
Set-Cell $ws 7 3 43945 2  # 2020-04-24
Set-Cell $ws 7 4 43957 1  # 2020-05-06
Set-Cell $ws 8 1 "Test models with mixed alloy data" 0
Set-Cell $ws 8 2 "Total data" 0
Set-Cell $ws 8 3 43964 2  # 2020-05-13
Set-Cell $ws 8 4 43971 1  # 2020-05-20
Set-Cell $ws 9 1 "Data Visualization: PCA, KNN, MDS" 3
Set-Cell $ws 9 2 "Total data" 0
Set-Cell $ws 9 3 43950 2  # 2020-04-29
Set-Cell $ws 9 4 43957 1  # 2020-05-06
Set-Cell $ws 10 1 "Expanded dataset (12k points) exploration" 3
Set-Cell $ws 10 2 "Expanded data" 0
Set-Cell $ws 10 3 43950 2  # 2020-04-29
Set-Cell $ws 10 4 43964 1  # 2020-05-13
Set-Cell $ws 11 1 "Brainstorm more descriptors" 0
Set-Cell $ws 11 2 "Lit Review" 2
Set-Cell $ws 11 3 43966 2  # 2020-05-15
Set-Cell $ws 11 4 43973 1  # 2020-05-22
Set-Cell $ws 12 1 "Expected improvement: informed selection of new training data" 3
Set-Cell $ws 12 2 "Bonus" 2
Set-Cell $ws 12 3 43971 2  # 2020-05-20
Set-Cell $ws 12 4 43979 1  # 2020-05-28
Set-Cell $ws 13 1 "Have model completed" 0
Set-Cell $ws 13 2 "Overall" 0
Set-Cell $ws 13 3 43987 2  # 2020-06-05
Set-Cell $ws 13 4 44001 2  # 2020-06-19

# Match formatting of the source sheets: column widths sized to content.
$ws.Columns.Item(1).ColumnWidth = 57
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 8.9
$ws.Columns.Item(4).ColumnWidth = 9.2

# Restore the selection on the previously-active "200422" sheet to the
# full used range, and drop its tab-selected state (select first so this
# doesn't leave it as the final active sheet).
$ws2.Range("A1:D11").Select()

# Make the new "200501" sheet the active tab with the same selection
# state recorded in the source file.
$ws.Range("I13").Select()

$wb.Save()

